# The sheet holds a data series in columns A/B (rows 2-101), with the very
# last row also carrying an (empty) marker cell in column C that denotes the
# end of the series. This edit appends 8 more data points to the series
# (rows 102-109) and moves that end-of-series marker from C101 down to C109.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the trailing marker cell from the old last row (101) to the new
# last row (109), then clear it from its old location.
$ws.Range("C101").Copy($ws.Range("C109"))
$ws.Range("C101").ClearContents()

# New rows 102-108 repeat the same A/B values.
for ($r = 102; $r -le 108; $r++) {
    $ws.Cells.Item($r, 1).Value = 0
    $ws.Cells.Item($r, 2).Value = 7.171333983999999
}

# The final row (109) has a distinct B value.
$ws.Cells.Item(109, 1).Value = 0
$ws.Cells.Item(109, 2).Value = 7.55965918
